$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.065647602081299
$ws.Range("B1").Value = 6.307534694671631
$ws.Range("C1").Value = 6.501516819000244
$ws.Range("D1").Value = 6.925323486328125
$ws.Range("E1").Value = 5.019288539886475
